$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D values that look like plain numbers must be forced to Text
# so Excel does not reinterpret/reformat them (e.g. trailing zeros).

$ws.Range("D2").Value = '28.108.07'
$ws.Range("E2").Value = '  +2.30%  '

$ws.Range("D3").Value = '1.909.95'
$ws.Range("E3").Value = '  +2.00%  '

$ws.Range("E4").Value = '  -1.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.73'
$ws.Range("E5").Value = '  +0.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  -1.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4829'
$ws.Range("E7").Value = '  +1.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3818'
$ws.Range("E8").Value = '  +1.29%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07356'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9359'
$ws.Range("E10").Value = '  -0.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.82'
$ws.Range("E11").Value = '  +0.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07790'
$ws.Range("E12").Value = '  -0.68%  '

$ws.Range("D13").Value = '1.910.97'
$ws.Range("E13").Value = '  +1.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.506'
$ws.Range("E14").Value = '  +1.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.624'
$ws.Range("E15").Value = '  +0.54%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.46'
$ws.Range("E16").Value = '  +0.60%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.004'
$ws.Range("E17").Value = '  -1.22%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008821'
$ws.Range("E18").Value = '  -1.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.003'
$ws.Range("E19").Value = '  -1.00%  '

$ws.Range("D20").Value = '28.119.99'
$ws.Range("E20").Value = '  +2.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.84'
$ws.Range("E21").Value = '  -0.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.142'

$ws.Range("D23").Value = '2.153.93'
$ws.Range("E23").Value = '  +0.94%  '

$ws.Range("E24").Value = '  +1.49%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.81'
$ws.Range("E25").Value = '  +1.90%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.926'
$ws.Range("E26").Value = '  -2.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.56'
$ws.Range("E27").Value = '  +0.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.108'
$ws.Range("E28").Value = '  +4.43%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '116.44'
$ws.Range("E29").Value = '  +0.45%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.950'
$ws.Range("E30").Value = '  -0.92%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08917'
$ws.Range("E31").Value = '  -0.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.341'
$ws.Range("E32").Value = '  -0.07%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.253'
$ws.Range("E33").Value = '  +3.03%  '

$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7687'
$ws.Range("E34").Value = '  +2.26%  '

$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.685'
$ws.Range("E35").Value = '  +1.42%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.609'
$ws.Range("E36").Value = '  -2.91%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02058'
$ws.Range("E37").Value = '  -0.45%  '

$ws.Range("E38").Value = '  -1.37%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05302'
$ws.Range("E39").Value = '  +0.00%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5484'
$ws.Range("E40").Value = '  +2.54%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.974'
$ws.Range("E41").Value = '  -1.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.026'
$ws.Range("E42").Value = '  -0.85%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1523'
$ws.Range("E43").Value = '  -0.22%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.464'
$ws.Range("E44").Value = '  +0.48%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.69'
$ws.Range("E45").Value = '  +0.73%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4832'
$ws.Range("E46").Value = '  -0.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '107.20'
$ws.Range("E47").Value = '  +4.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.003'
$ws.Range("E48").Value = '  -1.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.658'
$ws.Range("E49").Value = '  -0.24%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '68.24'
$ws.Range("E50").Value = '  +1.30%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06105'
$ws.Range("E51").Value = '  +0.19%  '
